$d = $word.ActiveDocument

# Helper: replace a run's text content while preserving sibling runs (e.g. an
# adjacent empty <w:r/>) in the same paragraph. Plain Find/Replace can cause
# the interop runtime to drop a purely-empty sibling run when the matched run
# itself carries no run-level formatting, so for those cases we locate the
# match, then re-insert just a plain run with the new text into that exact
# span via InsertXML (which, applied to a fresh Range object, replaces only
# the targeted span and leaves the rest of the paragraph/document intact).
function Set-PlainRunText($doc, $oldText, $newText) {
    $searchRng = $doc.Content
    $found = $searchRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $oldText"
    }
    $target = $doc.Range($searchRng.Start, $searchRng.End)
    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

# 1. Title / heading text (appears twice: H1 heading and bold run near the
#    bottom). Both occurrences share identical text, so a single global
#    Find/Replace updates them together; this keeps their bold run's
#    formatting intact.
$d.Content.Find.Execute(
    "Play Kings of Gold for Free - iSoftBet's New Egyptian Themed Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Kings of Gold Slot Free | Review and Gameplay", 2)

# 2. "What we like" bullet list items
Set-PlainRunText $d "Hold & Win bonus game for extra rewards" "Medium rewards up to 1,121x the bet"
Set-PlainRunText $d "Free spins feature triggered by the Pyramid scatter symbol" "Hold & Win bonus game"
Set-PlainRunText $d "Sticky gold coin symbols in Respins feature add winning opportunities" "Free spins feature"
Set-PlainRunText $d "Opulent Egyptian theme and background inspired design" "Opulent and luxurious design"

# 3. "What we don't like" bullet list items
Set-PlainRunText $d "Below-average RTP rate at 95.14%" "Below-average RTP of 95.14%"
Set-PlainRunText $d "Limited maximum win of 1,121x total bet" "Limited number of free spins symbols"

# 4. Meta description (italic run near the end)
$d.Content.Find.Execute(
    "Experience the opulence of ancient Egypt with Kings of Gold slot game's features such as Hold & Win bonus, free spins and RTP of 95.14%. Play free now.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Kings of Gold slot game and play for free to experience ancient Egypt themed fun.", 2)
